$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2-3: cardholder name / account number update
$ws.Range("C2").Value = "Hartmut"

# B3 holds a long numeric-looking card/account number that must stay text.
# Assign with a leading apostrophe to force text, then re-paste the original
# cell's formatting (a plain unformatted numeric-style cell) so the style
# index is not perturbed by the quote-prefix flag.
$ws.Range("B3").Value = "'2570314725427075"
$ws.Range("D3").Copy()
$ws.Range("B3").PasteSpecial(-4122)

$ws.Range("C3").Value = "Mohaupt"

# Row 5: opening balance date
$ws.Range("D5").Value = "KONTOSTAND AM 17.03.2025"

# Row 6
$ws.Range("B6").Value = "21.03."
$ws.Range("C6").Value = "22.03."
$ws.Range("D6").Value = "PAYPAL LNCUTB"
$ws.Range("E6").Value = "97,01-"

# Row 7
$ws.Range("B7").Value = "25.03."
$ws.Range("C7").Value = "26.03."
$ws.Range("D7").Value = "BEITRAG Allianz SE K-75866988"
$ws.Range("E7").Value = "56,80-"

# Row 8
$ws.Range("B8").Value = "29.03."
$ws.Range("C8").Value = "30.03."
$ws.Range("D8").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 85160029"
$ws.Range("E8").Value = "86,59-"

# Row 9
$ws.Range("B9").Value = "02.04."
$ws.Range("C9").Value = "03.04."
$ws.Range("D9").Value = "KARTENZ./02.04 ALDI SUED RO"
$ws.Range("E9").Value = "98,71-"

# Row 10: transaction removed entirely (was 01.01./02.01./KARTENZ... row) -
# clear the contents, leaving B10/C10/D10 blank and E10 blank too but with a
# slightly different (right+vcenter+wrap) alignment than its previous
# right-only alignment.
$ws.Range("B10").Value = $null
$ws.Range("C10").Value = $null
$ws.Range("D10").Value = $null
$ws.Range("E10").Value = $null
$ws.Range("E10").HorizontalAlignment = -4152
$ws.Range("E10").VerticalAlignment = -4108
$ws.Range("E10").WrapText = $true

# Row 11: transaction removed entirely (was 04.01./05.01./ABSCHLAG STROM...)
$ws.Range("B11").Value = $null
$ws.Range("C11").Value = $null
$ws.Range("D11").Value = $null
$ws.Range("E11").Value = $null
$ws.Range("E11").HorizontalAlignment = -4152
$ws.Range("E11").VerticalAlignment = -4108
$ws.Range("E11").WrapText = $true

# Row 12: closing balance date/amount
$ws.Range("D12").Value = "KONTOSTAND AM 04.04.2025"
$ws.Range("E12").Value = "339,11-"

# Row 13: next billing date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 13.04.2025"
